$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.032.89"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.695.83"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  +1.90%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "611.67"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.47"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +1.02%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.590"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E9").Value = "  +3.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.03"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +4.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.403"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -1.37%  "
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("E13").Value = "  +9.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "30.12"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +2.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.179.97"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.898.32"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.692.48"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +1.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.77"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.90"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.79"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +5.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "359.21"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.40"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +3.05%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("E24").Value = "  +17.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.98"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +5.67%  "
$ws.Range("E26").Value = "  -3.82%  "
$ws.Range("E27").Value = "  +0.40%  "
$ws.Range("E28").Value = "  +3.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.30"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.21"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "533.75"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -2.33%  "
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.70"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +4.45%  "
$ws.Range("E35").Value = "  -1.60%  "
$ws.Range("E36").Value = "  +0.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.77"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +0.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "162.06"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("E39").Value = "  -1.42%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.57"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "168.07"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.17"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0635"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +2.63%  "
$ws.Range("E46").Value = "  +2.39%  "
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0268"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +1.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.658"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("E50").Value = "  +6.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0998"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +1.39%  "
